# Scheduled market-data refresh: updates the pricing/profit columns
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ -> columns H-N)
# for a handful of leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 208.82608
$ws.Range("I33").Value = 217.40909
$ws.Range("K33").Value = 217.40909
$ws.Range("M33").Value = 11.59091000000001
$ws.Range("H43").Value = 4256.2856
$ws.Range("I43").Value = 4256.2856
$ws.Range("K43").Value = 4256.2856
$ws.Range("M43").Value = -4187.2856
$ws.Range("H74").Value = 4643.75
$ws.Range("I74").Value = 4420
$ws.Range("K74").Value = 4420
$ws.Range("M74").Value = -3484
$ws.Range("H77").Value = 4643.75
$ws.Range("I77").Value = 4420
$ws.Range("K77").Value = 22100
$ws.Range("M77").Value = -17420
$ws.Range("H116").Value = 11381.704
$ws.Range("I116").Value = 6511.4287
$ws.Range("K116").Value = 6511.4287
$ws.Range("M116").Value = -3069.4287
$ws.Range("H125").Value = 8205.833000000001
$ws.Range("I125").Value = 7449.5
$ws.Range("J125").Value = 8584
$ws.Range("K125").Value = 67045.5
$ws.Range("L125").Value = 77256
$ws.Range("M125").Value = -64585.5
$ws.Range("N125").Value = -82176
$ws.Range("H135").Value = 2520.762
$ws.Range("I135").Value = 733.4737
$ws.Range("J135").Value = 19500
$ws.Range("K135").Value = 6601.2633
$ws.Range("L135").Value = 175500
$ws.Range("M135").Value = -4066.2633
$ws.Range("N135").Value = -180570
$ws.Range("H138").Value = 3997.0544
$ws.Range("I138").Value = 2870.2727
$ws.Range("K138").Value = 8610.8181
$ws.Range("M138").Value = -3470.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3844.0823
$ws.Range("I32").Value = 3844.0823
$ws.Range("K32").Value = 3844.0823
$ws.Range("M32").Value = -3557.0823
$ws.Range("H37").Value = 38575
$ws.Range("I37").Value = 30034
$ws.Range("J37").Value = 39998.5
$ws.Range("K37").Value = 30034
$ws.Range("L37").Value = 39998.5
$ws.Range("M37").Value = -29761
$ws.Range("N37").Value = -40544.5
$ws.Range("H61").Value = 6182827.5
$ws.Range("I61").Value = 7413349.5
$ws.Range("K61").Value = 7413349.5
$ws.Range("M61").Value = -7413137.5
$ws.Range("H74").Value = 3222.4167
$ws.Range("I74").Value = 3218.476
$ws.Range("K74").Value = 3218.476
$ws.Range("M74").Value = -2344.476
$ws.Range("H77").Value = 3222.4167
$ws.Range("I77").Value = 3218.476
$ws.Range("K77").Value = 16092.38
$ws.Range("M77").Value = -11724.38
$ws.Range("H122").Value = 3588
$ws.Range("I122").Value = 3548.5
$ws.Range("K122").Value = 10645.5
$ws.Range("M122").Value = -8195.5
$ws.Range("H136").Value = 6182827.5
$ws.Range("I136").Value = 7413349.5
$ws.Range("K136").Value = 22240048.5
$ws.Range("M136").Value = -22237498.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 14666
$ws.Range("I107").Value = 19999
$ws.Range("K107").Value = 19999
$ws.Range("M107").Value = -18079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 923.1539
$ws.Range("J22").Value = 1736.25
$ws.Range("L22").Value = 1736.25
$ws.Range("N22").Value = -2436.25
$ws.Range("H58").Value = 2016.6216
$ws.Range("I58").Value = 1885.6666
$ws.Range("K58").Value = 1885.6666
$ws.Range("M58").Value = -1682.6666
$ws.Range("H86").Value = 18618.75
$ws.Range("I86").Value = 21400
$ws.Range("J86").Value = 13983.333
$ws.Range("K86").Value = 21400
$ws.Range("L86").Value = 13983.333
$ws.Range("M86").Value = -20277
$ws.Range("N86").Value = -16229.333
$ws.Range("H89").Value = 18618.75
$ws.Range("I89").Value = 21400
$ws.Range("J89").Value = 13983.333
$ws.Range("K89").Value = 107000
$ws.Range("L89").Value = 69916.66500000001
$ws.Range("M89").Value = -101384
$ws.Range("N89").Value = -81148.66500000001
$ws.Range("H122").Value = 3470.5881
$ws.Range("I122").Value = 3117.75
$ws.Range("K122").Value = 9353.25
$ws.Range("M122").Value = -6903.25
$ws.Range("H134").Value = 1683.7142
$ws.Range("I134").Value = 1781.5333
$ws.Range("K134").Value = 5344.5999
$ws.Range("M134").Value = -2809.5999
$ws.Range("H136").Value = 2016.6216
$ws.Range("I136").Value = 1885.6666
$ws.Range("K136").Value = 5656.9998
$ws.Range("M136").Value = -3106.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 16072
$ws.Range("J62").Value = 16072
$ws.Range("L62").Value = 48216
$ws.Range("N62").Value = -49588
$ws.Range("H65").Value = 16072
$ws.Range("J65").Value = 16072
$ws.Range("L65").Value = 144648
$ws.Range("N65").Value = -151512
$ws.Range("H82").Value = 27674.5
$ws.Range("J82").Value = 27674.5
$ws.Range("L82").Value = 83023.5
$ws.Range("N82").Value = -83835.5
$ws.Range("H85").Value = 27674.5
$ws.Range("J85").Value = 27674.5
$ws.Range("L85").Value = 83023.5
$ws.Range("N85").Value = -85831.5
$ws.Range("H131").Value = 3874.2
$ws.Range("I131").Value = 2038.8334
$ws.Range("J131").Value = 8593.714
$ws.Range("K131").Value = 6116.5002
$ws.Range("L131").Value = 25781.142
$ws.Range("M131").Value = -1076.5002
$ws.Range("N131").Value = -35861.142
$ws.Range("H132").Value = 2209.6538
$ws.Range("I132").Value = 1912.625
$ws.Range("J132").Value = 2684.9
$ws.Range("K132").Value = 17213.625
$ws.Range("L132").Value = 24164.1
$ws.Range("M132").Value = -14683.625
$ws.Range("N132").Value = -29224.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2271.0625
$ws.Range("J102").Value = 2644
$ws.Range("L102").Value = 2644
$ws.Range("N102").Value = -5888
$ws.Range("H113").Value = 1428031.1
$ws.Range("I113").Value = 3821.2856
$ws.Range("K113").Value = 3821.2856
$ws.Range("M113").Value = -1651.2856
$ws.Range("H122").Value = 6081
$ws.Range("I122").Value = 3896.75
$ws.Range("J122").Value = 8993.333000000001
$ws.Range("K122").Value = 11690.25
$ws.Range("L122").Value = 26979.999
$ws.Range("M122").Value = -9240.25
$ws.Range("N122").Value = -31879.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9750.049999999999
$ws.Range("I7").Value = 8002.4546
$ws.Range("J7").Value = 11886
$ws.Range("K7").Value = 8002.4546
$ws.Range("L7").Value = 11886
$ws.Range("M7").Value = -7890.4546
$ws.Range("N7").Value = -12110
$ws.Range("H22").Value = 23479.8
$ws.Range("I22").Value = 38599.668
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 38599.668
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -38304.668
$ws.Range("N22").Value = -1390
$ws.Range("H27").Value = 23479.8
$ws.Range("I27").Value = 38599.668
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 38599.668
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -38492.668
$ws.Range("N27").Value = -1014
$ws.Range("H40").Value = 8499.625
$ws.Range("J40").Value = 6949.5
$ws.Range("L40").Value = 6949.5
$ws.Range("N40").Value = -7221.5
$ws.Range("H126").Value = 9750.049999999999
$ws.Range("I126").Value = 8002.4546
$ws.Range("J126").Value = 11886
$ws.Range("K126").Value = 24007.3638
$ws.Range("L126").Value = 35658
$ws.Range("M126").Value = -21537.3638
$ws.Range("N126").Value = -40598
$ws.Range("H132").Value = 3268.658
$ws.Range("I132").Value = 2270.3684
$ws.Range("K132").Value = 6811.1052
$ws.Range("M132").Value = -4281.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5154.5
$ws.Range("J126").Value = 4397.3335
$ws.Range("L126").Value = 13192.0005
$ws.Range("N126").Value = -18132.0005
$ws.Range("H132").Value = 239362
$ws.Range("I132").Value = 859.84375
$ws.Range("K132").Value = 2579.53125
$ws.Range("M132").Value = -49.53125
